# Add a clarifying note textbox to slide 3 ("Summary of data conversion process"),
# matching the look & feel of the existing note textbox on that slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The existing "Note: ..." callout box is shape index 4 (id=2, "TextBox 1").
# Duplicate it so the new box inherits identical formatting (fill, body, font).
$template = $s.Shapes.Item(4)
$newShape = $template.Duplicate()

$newShape.Name = "TextBox 2"

# Position/size taken from the target OOXML (EMU / 12700 = points), using
# decimal values that round-trip exactly back to the target EMU values.
$newShape.Left = 136.7693
$newShape.Top = 386.1915
$newShape.Width = 552.71543
$newShape.Height = 29.0813

$newShape.TextFrame.TextRange.Text = "Note: the term " + [char]8220 + "client" + [char]8221 + " in this slide deck refers to an Enterprise using gist."
